$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (new content sourced from original row 4)
$ws.Range("A3").Value = 131039759
$ws.Range("B3").Value = 91829
$ws.Range("E3").Value = 5442
$ws.Range("F3").Value = 'Tallticka'
$ws.Range("G3").Value = 'Porodaedalea pini'
$ws.Range("H3").Value = '(Brot.) Murrill'
$ws.Range("P3").Value = 'Gotvad, Dlr'
$ws.Range("Q3").Value = 479059
$ws.Range("R3").Value = 6792254
$ws.Range("S3").Value = 10
$ws.Range("AC3").ClearContents()

# Row 4 (new content sourced from original row 3)
$ws.Range("A4").Value = 131041641
$ws.Range("B4").Value = 79243
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = 'Garnlav'
$ws.Range("G4").Value = 'Alectoria sarmentosa'
$ws.Range("H4").Value = '(Ach.) Ach.'
$ws.Range("P4").Value = 'Tandbergsvasseln, Dlr'
$ws.Range("Q4").Value = 479078
$ws.Range("R4").Value = 6791615
$ws.Range("S4").Value = 50
$ws.Range("AC4").Value = 'Rikligt i en radie av ca 50 meter, synfältet'

# Row 7 (new content sourced from original row 8)
$ws.Range("A7").Value = 131042226
$ws.Range("B7").Value = 79243
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = 'Garnlav'
$ws.Range("G7").Value = 'Alectoria sarmentosa'
$ws.Range("H7").Value = '(Ach.) Ach.'
$ws.Range("M7").ClearContents()
$ws.Range("Q7").Value = 479114
$ws.Range("R7").Value = 6792438
$ws.Range("S7").Value = 50
$ws.Range("Z7").Value = '16:15'
$ws.Range("AB7").Value = '16:15'
$ws.Range("AC7").Value = 'Rikligt till måttligt i en radie av ca 50 meter,synfältet'

# Row 8 (new content sourced from original row 7)
$ws.Range("A8").Value = 131038653
$ws.Range("B8").Value = 57881
$ws.Range("E8").Value = 100049
$ws.Range("F8").Value = 'Spillkråka'
$ws.Range("G8").Value = 'Dryocopus martius'
$ws.Range("H8").Value = '(Linnaeus, 1758)'
$ws.Range("M8").Value = 'äldre spår'
$ws.Range("Q8").Value = 479094
$ws.Range("R8").Value = 6792753
$ws.Range("S8").Value = 10
$ws.Range("Z8").Value = '12:08'
$ws.Range("AB8").Value = '12:08'
$ws.Range("AC8").ClearContents()

# Row 9 (new content sourced from original row 10)
$ws.Range("A9").Value = 131040374
$ws.Range("B9").Value = 79001
$ws.Range("E9").Value = 228912
$ws.Range("F9").Value = 'Mörk kolflarnlav'
$ws.Range("G9").Value = 'Carbonicola myrmecina'
$ws.Range("H9").Value = '(Ach.) Bendiksby & Timdal'
$ws.Range("Q9").Value = 479088
$ws.Range("R9").Value = 6792211

# Row 10 (new content sourced from original row 9)
$ws.Range("A10").Value = 131039523
$ws.Range("B10").Value = 79243
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = 'Garnlav'
$ws.Range("G10").Value = 'Alectoria sarmentosa'
$ws.Range("H10").Value = '(Ach.) Ach.'
$ws.Range("Q10").Value = 479079
$ws.Range("R10").Value = 6792517

# Row 12 (new content sourced from original row 13)
$ws.Range("A12").Value = 131039119
$ws.Range("B12").Value = 79243
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = 'Garnlav'
$ws.Range("G12").Value = 'Alectoria sarmentosa'
$ws.Range("H12").Value = '(Ach.) Ach.'
$ws.Range("Q12").Value = 479105
$ws.Range("R12").Value = 6792638

# Row 13 (new content sourced from original row 12)
$ws.Range("A13").Value = 131040483
$ws.Range("B13").Value = 78646
$ws.Range("E13").Value = 6437
$ws.Range("F13").Value = 'Blanksvart spiklav'
$ws.Range("G13").Value = 'Calicium denigratum'
$ws.Range("H13").Value = '(Vain.) Tibell'
$ws.Range("Q13").Value = 479088
$ws.Range("R13").Value = 6792211

# Row 18 (new content sourced from original row 19)
$ws.Range("A18").Value = 131041965
$ws.Range("B18").Value = 57884
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = 'Tretåig hackspett'
$ws.Range("G18").Value = 'Picoides tridactylus'
$ws.Range("H18").Value = '(Linnaeus, 1758)'
$ws.Range("M18").Value = 'färska spår'
$ws.Range("P18").Value = 'Tandbergsvasseln, Dlr'
$ws.Range("Q18").Value = 479096
$ws.Range("R18").Value = 6792085
$ws.Range("S18").Value = 10
$ws.Range("AC18").ClearContents()

# Row 19 (new content sourced from original row 20)
$ws.Range("A19").Value = 131039828
$ws.Range("M19").Value = 'bobygge'
$ws.Range("P19").Value = 'Gotvad, Dlr'
$ws.Range("Q19").Value = 479059
$ws.Range("R19").Value = 6792254
$ws.Range("AE19").Value = $true

# Row 20 (new content sourced from original row 18)
$ws.Range("A20").Value = 131039579
$ws.Range("B20").Value = 79243
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = 'Garnlav'
$ws.Range("G20").Value = 'Alectoria sarmentosa'
$ws.Range("H20").Value = '(Ach.) Ach.'
$ws.Range("M20").ClearContents()
$ws.Range("Q20").Value = 479079
$ws.Range("R20").Value = 6792475
$ws.Range("S20").Value = 50
$ws.Range("AC20").Value = 'Rikligt till måttligt i en radie av ca 50 meter, synfältet'
$ws.Range("AE20").Value = $false

# Row 26 (new content sourced from original row 27)
$ws.Range("A26").Value = 131039763
$ws.Range("Q26").Value = 479059
$ws.Range("R26").Value = 6792254
$ws.Range("S26").Value = 10
$ws.Range("AC26").ClearContents()

# Row 27 (new content sourced from original row 26)
$ws.Range("A27").Value = 131039672
$ws.Range("Q27").Value = 479066
$ws.Range("R27").Value = 6792326
$ws.Range("S27").Value = 50
$ws.Range("AC27").Value = 'Rikligt till måttligt i en radie av ca 50 meter, synfältet'
